$d = $word.ActiveDocument

# --- Rewrite the "theory strengthening" sentences in the first body
#     paragraph under "How does it work?" ---

$d.Content.Find.Execute(
    "Different kinds of entries will hold higher merit than other kinds.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Different types of entries will hold higher merit than others.",
    2) | Out-Null

$d.Content.Find.Execute(
    "The theory of an entry is strengthened by the number of times reproduced, peer review, and child tests within and expanding the domain of the hypothesis" + [char]0x2019 + " predictions.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If a submission fails to disprove its parent, then it strengthens the theory the parent" + [char]0x2019 + "s hypothesis.",
    2) | Out-Null

# --- Add a blank paragraph after the "Who can work on this project?"
#     answer (the GitHub paragraph), right before the trailing bookmark
#     paragraph. ---

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*psychological effect.*") {
        $target = $i
    }
}
if ($target -ne $null) {
    $bookmarkPara = $d.Paragraphs($target + 1)
    $bookmarkPara.Range.InsertParagraphBefore()
}
